$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 43
$ws.Range("B4").Value = "Mazeu"
$ws.Range("C4").Value = "maze@gmail.com"
$ws.Range("E4").Value = "2023-06-12T21:00:44.000000Z"
$ws.Range("F4").Value = "2023-06-12T21:00:44.000000Z"
